$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.149.99'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '3.813.54'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '707.15'
$ws.Range("E5").Value = '  +0.92%  '
$ws.Range("D6").Value = '171.74'
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("D7").Value = '3.811.98'
$ws.Range("E7").Value = '  -1.28%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").Value = '7.73'
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  -2.38%  '
$ws.Range("D14").Value = '35.91'
$ws.Range("E14").Value = '  -1.21%  '
$ws.Range("D15").Value = '4.454.60'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").Value = '3.794.46'
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").Value = '71.099.16'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").Value = '17.44'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '503.67'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Value = '10.72'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '84.35'
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("E25").Value = '  -3.34%  '
$ws.Range("D26").Value = '3.963.21'
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").Value = '12.07'
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("D28").Value = '10.39'
$ws.Range("E28").Value = '  -2.70%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  -4.32%  '
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("E32").Value = '  -0.79%  '
$ws.Range("D33").Value = '7.36'
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("D34").Value = '29.06'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("D35").Value = '0.174'
$ws.Range("E35").Value = '  -4.47%  '
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").Value = '3.777.09'
$ws.Range("E37").Value = '  -1.09%  '
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("E39").Value = '  -2.58%  '
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("E41").Value = '  -2.99%  '
$ws.Range("D42").Value = '5.94'
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("D43").Value = '3.27'
$ws.Range("E43").Value = '  -4.78%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '166.91'
$ws.Range("E46").Value = '  +1.94%  '
$ws.Range("D47").Value = '0.000314'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").Value = '49.13'
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("D49").Value = '421.06'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("D50").Value = '8.62'
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = '1.16'
$ws.Range("E51").Value = '  +3.36%  '
